# ModificarCliente.xlsx - "update entregable 1, 2"
#
# The underlying change is a data correction on row 2 of the
# "ModificarCliente" sheet:
#   - E Civil (column H)  : "Casado" -> "CASADO"
#   - Fecha    (column M)  : "3 jul. 2023, 14:55:25" -> "14 jul. 2023, 09:44:27"
#
# (All of the shared-string index churn visible in the raw XML diff is a
# side effect of these two textual edits causing Excel to re-pack the
# shared string table; the actual displayed values for every other cell
# are unchanged.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Estado Civil" value for the row-2 record. Setting .Value
# directly resets the cell's number-format/fill style, so re-apply the
# original formatting (shared with its row-2 neighbours, e.g. I2) by
# copying just the formats back onto the cell afterwards.
$ws.Range("H2").Value = "CASADO"
$ws.Range("I2").Copy()
$ws.Range("H2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the "Fecha" (timestamp) value for the row-2 record.
$ws.Range("M2").Value = "14 jul. 2023, 09:44:27"

# Mirror the saved cursor/selection position recorded in the workbook.
$ws.Range("M4").Select()
